$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(3, 6).Value = 238
$ws.Cells.Item(5, 6).Value = 1759
$ws.Cells.Item(6, 6).Value = 652
$ws.Cells.Item(8, 6).Value = 448
$ws.Cells.Item(9, 6).Value = 4305
$ws.Cells.Item(11, 6).Value = 447
$ws.Cells.Item(17, 6).Value = 2946
$ws.Cells.Item(18, 6).Value = 1777
$ws.Cells.Item(21, 6).Value = 161
$ws.Cells.Item(23, 6).Value = 913
$ws.Cells.Item(24, 6).Value = 288
$ws.Cells.Item(26, 6).Value = 2244
$ws.Cells.Item(28, 6).Value = 2311
$ws.Cells.Item(30, 6).Value = 683
$ws.Cells.Item(31, 6).Value = 497
$ws.Cells.Item(33, 6).Value = 875
$ws.Cells.Item(34, 6).Value = 400
$ws.Cells.Item(35, 6).Value = 1067
$ws.Cells.Item(36, 6).Value = 879
$ws.Cells.Item(37, 6).Value = 1144
$ws.Cells.Item(38, 6).Value = 7
$ws.Cells.Item(39, 6).Value = 323
$ws.Cells.Item(40, 6).Value = 503
$ws.Cells.Item(42, 6).Value = 270
$ws.Cells.Item(43, 6).Value = 3452

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(3, 6).Value = 238
$ws.Cells.Item(6, 6).Value = 1759
$ws.Cells.Item(7, 6).Value = 651
$ws.Cells.Item(9, 6).Value = 448
$ws.Cells.Item(10, 6).Value = 4305
$ws.Cells.Item(17, 6).Value = 2946
$ws.Cells.Item(19, 6).Value = 1777
$ws.Cells.Item(22, 6).Value = 161
$ws.Cells.Item(27, 6).Value = 913
$ws.Cells.Item(28, 6).Value = 288
$ws.Cells.Item(29, 6).Value = 2244
$ws.Cells.Item(33, 6).Value = 2311
$ws.Cells.Item(34, 6).Value = 683
$ws.Cells.Item(35, 6).Value = 497
$ws.Cells.Item(36, 6).Value = 875
$ws.Cells.Item(37, 6).Value = 1067
$ws.Cells.Item(38, 6).Value = 879
$ws.Cells.Item(39, 6).Value = 1144
$ws.Cells.Item(40, 6).Value = 323
$ws.Cells.Item(41, 6).Value = 503
$ws.Cells.Item(47, 6).Value = 270
$ws.Cells.Item(48, 6).Value = 3452
